$wb = $excel.ActiveWorkbook

# Duplicate the most recent daily sheet ("2025-08-19") to create the new
# "2025-08-20" sheet, placed right after it, mirroring the existing pattern.
$count = $wb.Worksheets.Count
$source = $wb.Worksheets.Item($count)
$source.Copy($null, $source)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "2025-08-20"

# New ranking data for 2025-08-20: rank, title, author, latest_episode
$data = @(
    @(1, "転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～", "zunta(作画) はらわたさいぞう(原作)", "第31話：完全なる死角③"),
    @(2, "時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―", "光永康則", "第６８話『施錠停止』③"),
    @(3, "地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。", "マツモトケンゴ", "第６３話　ダンスゲームの戦いが始まった（１）"),
    @(4, "勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが", "絶叫あいす。(漫画) 木の芽(原作) 希(キャラクター原案)", "第4話 後編"),
    @(5, "世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜", "戸賀 環 坂木持丸 riritto", "第51話①　呪われた家を探索してみた"),
    @(6, "バキ外伝 烈海王は異世界転生しても一向にかまわんッッ", "板垣恵介 猪原賽 陸井栄史", "第78話　海神(ポセイドン)"),
    @(7, "バキ外伝　ガイアとシコルスキー　～ときどきノムラ 二人だけど三人暮らし～", "板垣恵介 林たかあき", "第52話 洗練されたデザイン"),
    @(8, "お気楽領主の楽しい領地防衛 ～生産系魔術で名もなき村を最強の城塞都市に～", "青色まろ（漫画） 赤池宗（原作） 転（原作イラスト）", "第33話　観光案内"),
    @(9, "ありふれた職業で世界最強", "RoGa（漫画） 白米 良（原作） たかやKi（キャラクター原案）", "第84話　人間らしさ（前編）"),
    @(10, "婚約者に裏切られた錬金術師は、独立して『ざまぁ』します　コミック版", "漫画/すたひろ 原作/Y.A", "chapter68【36話①】"),
    @(11, "Lv２からチートだった元勇者候補のまったり異世界ライフ", "糸町秋音（漫画） 鬼ノ城ミヤ（原作） 片桐（キャラクター原案）", "第60話　ふたつの希望（前編）"),
    @(12, "女子高生の無駄づかい", "ビーノ(著者)", "第134話　かくせい"),
    @(13, "ひとりぼっちの異世界攻略", "びび（漫画） 五示正司（原作）", "第232話　ボコ地獄"),
    @(14, "Sランク冒険者である俺の娘たちは重度のファザコンでした", "しゅにち（漫画） 友橋かめつ（原作） 希望つばめ（原作イラスト）", "第47話　賢者vs大賢者"),
    @(15, "序盤で死ぬ最強のサブキャラに転生したので、ゲーム知識で無双する", "作画：マエD 原作：新人", "第6話(1)"),
    @(16, "勇者と魔物の運び屋さん", "大河原", "最終話 元の世界へ"),
    @(17, "新米オッサン冒険者、最強パーティに死ぬほど鍛えられて無敵になる", "漫画：荻野ケン 原作：岸馬きらく キャラクター原案：Tea", "第70話"),
    @(18, "江戸前エルフ", "樋口彰彦", "#119"),
    @(19, "元最強探索者のおじさん。美少女配信者を助けて大バズりしてしまった。", "かなたろー(原作) 草壁レイ(漫画)", "第4話　美少女、おじさんに「わからせ」られる。②"),
    @(20, "俺にトラウマを与えた女子達がチラチラ見てくるけど、残念ですが手遅れです", "いちたか（漫画） 御堂ユラギ（原作） 緜（キャラクター原案）", "第19話　合コンに参加（後編）"),
    @(21, "えろいことするために巨乳美少女奴隷を買ったはずが、お師匠さまと慕われて思った通りにいかなくなる話", "佐藤36(作画) 煮豆シューター(原作)", "第4話後半"),
    @(22, "最凶の支援職【話術士】である俺は世界最強クランを従える", "やもりちゃん（漫画） じゃき（原作）", "第54話　欺瞞と真相（前編）"),
    @(23, "最強の少年聖騎士、転生者を狩る", "作画：御塩 原作：宇奈木ユラ", "第7話(2)"),
    @(24, "みつばものがたり 呪いの少女と死の輪舞《ロンド》", "堤りん(漫画) 七沢またり(原作) EURA(キャラクター原案)", "第11話：勝利の美酒"),
    @(25, "最強勇者パーティーは愛が知りたい", "山田肌襦袢", "第29話「きみがきみであればいい」"),
    @(26, "うちの清楚系委員長がかつて中二病アイドルだったことを俺だけが知っている。", "三上こた こばやし少女 寝子空兄 ゆがー", "第1話　災禍の悪夢"),
    @(27, "魔法歌姫マジカルギンガ", "超銀河レコード", "第13話"),
    @(28, "四姉妹は夜をおまちかね", "保住圭(原作) Bcoca(作画)", "第35夜②：『特別』の形はひとつじゃなくて"),
    @(29, "死ぬ運命にある悪役令嬢の兄に転生したので、妹を育てて未来を変えたいと思います　～世界最強はオレだけど、世界最カワは妹に違いない～", "石見翔子(漫画） 泉里侑希（原作） タムラヨウ（キャラクター原案）", "第4話　はじめてのお出かけ（後編）"),
    @(30, "ラーメン大好き小泉さん", "鳴見なる", "21杯目 背脂"),
    @(31, "勇者パーティを追放された【スキルサポーター】、仲間のスキルを解放して最強に成り上がる", "作画：なかお 原作：前田氏", "第7話(2)"),
    @(32, "ネトゲの嫁が人気アイドルだった　～クール系の彼女は現実でも嫁のつもりでいる～", "もっつぉ（漫画） あボーン（原作） 館田ダン（キャラクター原案）", "第31話　黒い平原コラボカフェ"),
    @(33, "恋人のフリのフリ", "むねひろ(著者)", "第8話②"),
    @(34, "ゲーム　オブ　ファミリア-家族戦記-", "Ｄ．Ｐ(作画) 山口ミコト(原作)", "第74話③"),
    @(35, "現代ダンジョンライフの続きは異世界オープンワールドで！", "田中清久（漫画） しば犬部隊（原作） ひろせ（原作イラスト）", "第24話　お前の名前を――Ⅱ（後編）"),
    @(36, "勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～", "漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり", "第５１話　英雄を倒す器用貧乏（３）"),
    @(37, "初歩魔法しか使わない謎の老魔法使いが旅をする", "山代カゲツ(漫画) やまだのぼる(原作) にじまあるく(キャラクター原案)", "第6話①"),
    @(38, "【パクパクですわ】追放されたお嬢様の『モンスターを食べるほど強くなる』スキルは、１食で１レベルアップする前代未聞の最強スキルでした。３日で人類最強になりましたわ～！", "島知宏 音速炒飯 有都あらゆる", "第２３食　巨大ヘビモンスターさん、パクパクですわ！（１）"),
    @(39, "無能の中の無能王子　スキル【無能】を授かりましたが、周りの女性は【傾国】【傾城】【奸婦】【毒婦】【悪婦】【妖婦】とかです", "漫画/一夢 原作/福朗 キャラクター原案/菊池政治", "chapter11【6話②】"),
    @(40, "8歳から始める魔法学", "ゆうなぎ（漫画） 上野夕陽（原作） 乃希（キャラクター原案）", "第19話　小英雄譚"),
    @(41, "転生したら没落貴族だったので、【呪言】を極めて家族を救います", "作画：アマセケイ 原作：メソポ・たみあ", "第7話(2)"),
    @(42, "ちゃんと吸えない吸血鬼ちゃん", "二式恭介(著者)", "第103話：吸血鬼ちゃんと卒業旅行②"),
    @(43, "宇崎ちゃんは遊びたい！", "丈(著者)", "第126話"),
    @(44, "第七魔王子ジルバギアスの魔王傾国記", "野井ニトラ（漫画） 甘木智彬（原作） 輝竜 司（キャラクター原案）", "第17話　勇者の決断【急】（後編）"),
    @(45, "異世界で最強のスキルを生み出せたので、ひたすら無双することにしました。　～俺だけがスキルの数値を勝手に操作～", "漫画：星トマジロウ 原作：ヒゲ抜き地蔵 キャラクター原案：山椒魚", "第10話 ②"),
    @(46, "フシノカミ ～辺境から始める文明再生記～", "黒杞よるの（漫画） 雨川水海（原作） 大熊まい（キャラクター原案）", "第40話　積み重なる熱"),
    @(47, "おしかけ勇者嫁 勇者は放逐されたおっさんを追いかけ、スローライフを応援する　コミック版", "漫画/珠洲城くるみ 原作/日富美信吾 キャラクター原案/白井鋭利", "chapter33【16話②】"),
    @(48, "宮廷魔導師、追放される　～無能だと追い出された最巧の魔導師は、部下を引き連れて冒険者クランを始めるようです～", "きつね丸（漫画） しんこせい（原作） ろこ（キャラクター原案）", "第3話　誰かのための戦い（後編）"),
    @(49, "衛宮さんちの今日のごはん", "TAa(漫画) 只野まこと(料理監修) ＴＹＰＥ－ＭＯＯＮ(原作)", "第75話"),
    @(50, "異世界メイドの三ツ星グルメ ～現代ごはん作ったら王宮で大バズリしました～", "モリタ Ｕ４ nima", "第13話（２）　ドーナツの騎士様とお土産スイーツ（２）")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rdata = $data[$i]
    $newSheet.Cells.Item($row, 1).Value = $rdata[0]
    $newSheet.Cells.Item($row, 2).Value = $rdata[1]
    $newSheet.Cells.Item($row, 3).Value = $rdata[2]
    $newSheet.Cells.Item($row, 4).Value = $rdata[3]
}

# Restore the originally active sheet/tab (Sheet.Copy makes the new copy active)
$wb.Worksheets.Item(1).Activate()
